{"js": "// The document contains a single table of simple arithmetic prompts\n// (\"62-49=\", \"7+65=\", ...). The edit replaces the prompt text in each\n// cell, in table (reading) order, with a new prompt \u2014 formatting\n// (alignment, font, size) is left untouched; only the text changes.\nconst newValues = [\"33-27=\",\"18+53=\",\"91-26=\",\"16+76=\",\"43+38=\",\"9+62=\",\"26-19=\",\"5+89=\",\"18+3=\",\"44-26=\",\"93-59=\",\"53+9=\",\"29+58=\",\"92-75=\",\"43+49=\",\"17+68=\",\"70-13=\",\"70-11=\",\"39+56=\",\"17+37=\",\"14+9=\",\"56+8=\",\"32-24=\",\"83-44=\",\"40-36=\",\"19+5=\",\"17+24=\",\"91-77=\",\"8+6=\",\"83-37=\",\"27+44=\",\"29+36=\",\"81-44=\",\"9+19=\",\"52-5=\",\"51-2=\",\"75-47=\",\"80-31=\",\"8+7=\",\"66+16=\",\"24+59=\",\"43-35=\",\"25+7=\",\"41-23=\",\"79+16=\",\"83-27=\",\"50-22=\",\"56-28=\",\"81-59=\",\"92-13=\",\"77+5=\",\"29+29=\",\"93-18=\",\"29+15=\",\"44-18=\",\"41-12=\",\"79+5=\",\"66-57=\",\"6+59=\",\"31-7=\",\"89+6=\",\"66+8=\",\"74-59=\",\"29+44=\",\"3+9=\",\"20-18=\",\"95-57=\",\"37+35=\",\"74-67=\",\"81-78=\",\"39+36=\",\"53-19=\",\"9+57=\",\"18+44=\",\"62-13=\",\"19+57=\",\"26+17=\",\"37+18=\",\"45+6=\",\"90-16=\",\"48+14=\",\"38+45=\",\"20-13=\",\"80-38=\",\"28+37=\",\"44-29=\",\"7+36=\",\"55+29=\",\"4+19=\",\"29+26=\",\"19+74=\",\"15-8=\",\"19+78=\",\"16+27=\",\"43-35=\",\"77+7=\",\"49+44=\",\"52-23=\",\"92-4=\",\"26+8=\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nlet idx = 0;\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (const cell of cells.items) {\n    if (idx >= newValues.length) break;\n\n    const paragraphs = cell.body.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n\n    // Each cell holds exactly one paragraph with one run holding the\n    // prompt text; replace just that paragraph's text, preserving the\n    // run's formatting (font/size) and paragraph properties (alignment).\n    const para = paragraphs.items[0];\n    const range = para.getRange();\n    range.insertText(newValues[idx], Word.InsertLocation.replace);\n\n    idx++;\n  }\n  await context.sync();\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single table of simple arithmetic prompts\n# (\"62-49=\", \"7+65=\", ...). The edit replaces the prompt text in each\n# cell, in table (reading) order, with a new prompt -- formatting\n# (alignment, font, size) is left untouched; only the text changes.\n$newValues = @(\n  \"33-27=\",\"18+53=\",\"91-26=\",\"16+76=\",\"43+38=\",\n  \"9+62=\",\"26-19=\",\"5+89=\",\"18+3=\",\"44-26=\",\n  \"93-59=\",\"53+9=\",\"29+58=\",\"92-75=\",\"43+49=\",\n  \"17+68=\",\"70-13=\",\"70-11=\",\"39+56=\",\"17+37=\",\n  \"14+9=\",\"56+8=\",\"32-24=\",\"83-44=\",\"40-36=\",\n  \"19+5=\",\"17+24=\",\"91-77=\",\"8+6=\",\"83-37=\",\n  \"27+44=\",\"29+36=\",\"81-44=\",\"9+19=\",\"52-5=\",\n  \"51-2=\",\"75-47=\",\"80-31=\",\"8+7=\",\"66+16=\",\n  \"24+59=\",\"43-35=\",\"25+7=\",\"41-23=\",\"79+16=\",\n  \"83-27=\",\"50-22=\",\"56-28=\",\"81-59=\",\"92-13=\",\n  \"77+5=\",\"29+29=\",\"93-18=\",\"29+15=\",\"44-18=\",\n  \"41-12=\",\"79+5=\",\"66-57=\",\"6+59=\",\"31-7=\",\n  \"89+6=\",\"66+8=\",\"74-59=\",\"29+44=\",\"3+9=\",\n  \"20-18=\",\"95-57=\",\"37+35=\",\"74-67=\",\"81-78=\",\n  \"39+36=\",\"53-19=\",\"9+57=\",\"18+44=\",\"62-13=\",\n  \"19+57=\",\"26+17=\",\"37+18=\",\"45+6=\",\"90-16=\",\n  \"48+14=\",\"38+45=\",\"20-13=\",\"80-38=\",\"28+37=\",\n  \"44-29=\",\"7+36=\",\"55+29=\",\"4+19=\",\"29+26=\",\n  \"19+74=\",\"15-8=\",\"19+78=\",\"16+27=\",\"43-35=\",\n  \"77+7=\",\"49+44=\",\"52-23=\",\"92-4=\",\"26+8=\"\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    if ($idx -ge $newValues.Length) { break }\n    $cell = $tbl.Cell($r, $c)\n    # Assigning Range.Text replaces just the cell's text content while\n    # keeping the existing run/paragraph formatting (font, size, align).\n    $cell.Range.Text = $newValues[$idx]\n    $idx++\n  }\n}\n"}
